$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# --- Title paragraph: merge runs (cosmetic, text unchanged) ---
Replace-Text "Trabalho de Ionic 3 + Firebase Google + Angular 4" "Trabalho de Ionic 3 + Firebase Google + Angular 4"

# --- Professor name paragraph: merge runs (cosmetic, text unchanged) ---
Replace-Text "Prof. Dr. Daniel Facciolo Pires" "Prof. Dr. Daniel Facciolo Pires"

# --- Date: replace the literal date with "dia da prova" ---
Replace-Text "Data: 06/06/2018" "Data: dia da prova"

# --- Modo de Entrega: merge runs (cosmetic, text unchanged) ---
# ("Valor: 3,0 pts" is left untouched: its trailing proofErr sits right at the
#  end of the paragraph and a no-op replace there would strand an unpaired tag)
Replace-Text "Modo de Entrega: apresentar em aula" "Modo de Entrega: apresentar em aula"

# --- Utilizando Ionic 3... paragraph: merge runs (cosmetic, text unchanged) ---
Replace-Text "Utilizando Ionic 3, Firebase Google e Angular 4, faça uma aplicação móvel contendo as seguintes funcionalidades:" "Utilizando Ionic 3, Firebase Google e Angular 4, faça uma aplicação móvel contendo as seguintes funcionalidades:"

# --- Bullet list items: merge runs (cosmetic, text unchanged) ---
Replace-Text "Tela inicial com autenticação no sistema com email e senha; e botão para criar usuário" "Tela inicial com autenticação no sistema com email e senha; e botão para criar usuário"
Replace-Text "Tela para Registro de usuários com nome, username, telefone, CPF, email e senha" "Tela para Registro de usuários com nome, username, telefone, CPF, email e senha"
Replace-Text "Usuário não pode se registrar com email que já tenha sido registrado;" "Usuário não pode se registrar com email que já tenha sido registrado;"
Replace-Text "Usuário não pode se registrar com username que já tenha sido registrado;" "Usuário não pode se registrar com username que já tenha sido registrado;"

# --- Remove the "Tela para remover uma entidade;" bullet and the empty paragraph right after it ---
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.TrimEnd([char]13, [char]7) -eq "Tela para remover uma entidade;") {
        $target = $i
        break
    }
}
if ($target -ne $null) {
    $p1 = $d.Paragraphs($target)
    $p2 = $d.Paragraphs($target + 1)
    $r = $d.Range($p1.Range.Start, $p2.Range.End)
    $r.Delete()
}
